# "Cambios en las situaciones - Ivan Alejandro"
#
# Slide 7 ("Situación 2"):
#   - Reposition the "Rectángulo 3" text box (it moves up, out of the way
#     of the table that used to sit above it).
#   - Remove the "Tabla 4" comparison table entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Locate shapes by name so this is robust to any index shuffling.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Rectángulo 3") {
        $rect = $shape
    }
    if ($shape.Name -eq "Tabla 4") {
        $table = $shape
    }
}

# Move the rectangle text box to its new position.
# Target EMU offset: x=1324573, y=2210623 (EMU / 12700 = points).
$rect.Left = 1324573 / 12700
$rect.Top  = 2210623 / 12700

# Delete the comparison table shape.
$table.Delete()
